$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2, shifting existing data rows down by one
$ws.Rows("2:2").Insert()
$ws.Rows("2:2").ClearFormats()

# Populate the new row 2 with the new data record
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 44756
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 100112031
$ws.Range("G2").Value = "Poroto verde"
$ws.Range("H2").Value = "Magnum"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 28000
$ws.Range("L2").Value = 29000
$ws.Range("M2").Value = 28500
$ws.Range("N2").Value = "$/malla 25 kilos"
$ws.Range("O2").Value = "Perú"
$ws.Range("P2").Value = 1140
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"

# Restore the date cell's number format (style carried by column D elsewhere)
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
